$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "41.325.47"
Set-TextValue "E2" "  -3.03%  "

Set-TextValue "D3" "2.465.43"
Set-TextValue "E3" "  -2.16%  "

Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.05%  "

Set-TextValue "D5" "311.41"
Set-TextValue "E5" "  +0.86%  "

Set-TextValue "D6" "94.27"

Set-TextValue "E7" "  -3.16%  "

Set-TextValue "E8" "  +0.01%  "

Set-TextValue "D9" "0.502"
Set-TextValue "E9" "  -4.51%  "

Set-TextValue "D10" "33.61"
Set-TextValue "E10" "  -6.61%  "

Set-TextValue "E11" "  -2.77%  "

Set-TextValue "E12" "  -0.28%  "

Set-TextValue "D13" "6.95"
Set-TextValue "E13" "  -4.82%  "

Set-TextValue "D14" "2.845.31"
Set-TextValue "E14" "  -2.15%  "

Set-TextValue "D15" "2.492.70"
Set-TextValue "E15" "  -0.79%  "

Set-TextValue "E16" "  -6.26%  "

Set-TextValue "D17" "0.787"
Set-TextValue "E17" "  -2.33%  "

Set-TextValue "D18" "41.268.26"
Set-TextValue "E18" "  -3.08%  "

Set-TextValue "D19" "6.34"
Set-TextValue "E19" "  -5.49%  "

Set-TextValue "D20" "0.0₃0920"

Set-TextValue "D21" "11.48"
Set-TextValue "E21" "  -5.16%  "

Set-TextValue "D22" "67.85"
Set-TextValue "E22" "  -2.27%  "

Set-TextValue "D23" "235.97"
Set-TextValue "E23" "  -3.15%  "

Set-TextValue "D24" "2.78"
Set-TextValue "E24" "  -3.59%  "

Set-TextValue "D25" "1.93"
Set-TextValue "E25" "  -5.52%  "

Set-TextValue "E26" "  +0.10%  "

Set-TextValue "D27" "24.42"
Set-TextValue "E27" "  -5.89%  "

Set-TextValue "E28" "  -4.66%  "

Set-TextValue "D29" "9.70"
Set-TextValue "E29" "  -4.25%  "

Set-TextValue "D30" "36.13"
Set-TextValue "E30" "  -7.25%  "

Set-TextValue "D31" "152.94"
Set-TextValue "E31" "  -1.85%  "

Set-TextValue "D32" "5.55"
Set-TextValue "E32" "  -3.62%  "

Set-TextValue "D33" "2.60"
Set-TextValue "E33" "  -6.05%  "

Set-TextValue "D34" "2.60"
Set-TextValue "E34" "  -0.95%  "

Set-TextValue "D35" "0.0755"
Set-TextValue "E35" "  -4.10%  "

Set-TextValue "D36" "3.01"

Set-TextValue "E37" "  -6.09%  "

Set-TextValue "D38" "17.09"
Set-TextValue "E38" "  -6.00%  "

Set-TextValue "E39" "  -3.68%  "

Set-TextValue "D40" "0.103"
Set-TextValue "E40" "  -8.17%  "

Set-TextValue "D41" "4.25"
Set-TextValue "E41" "  +0.07%  "

Set-TextValue "D42" "21.23"
Set-TextValue "E42" "  -4.40%  "

Set-TextValue "E43" "  +0.10%  "

Set-TextValue "D44" "1.964.54"
Set-TextValue "E44" "  -0.99%  "

Set-TextValue "D45" "0.0285"
Set-TextValue "E45" "  -4.72%  "

Set-TextValue "D46" "3.07"
Set-TextValue "E46" "  -6.52%  "

Set-TextValue "D47" "8.66"
Set-TextValue "E47" "  -1.83%  "

Set-TextValue "D48" "69.91"
Set-TextValue "E48" "  -3.46%  "

Set-TextValue "D49" "76.09"
Set-TextValue "E49" "  -4.94%  "

Set-TextValue "D50" "97.26"
Set-TextValue "E50" "  -3.39%  "

Set-TextValue "E51" "  -5.90%  "
